$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "C6_L11_Ratioandproportions"

$ws.Range("B23").Select()
